# 2020.05.17 folder of web data mining 1 group
# Fill in the remaining group-member rows (6 and 7) of the reading-report
# summary table, add the running-total formulas for rows 5-9, and move the
# active selection to J7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (闻浩) ---------------------------------------------------------
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 25
$ws.Range("D6").Value = 3
$ws.Range("F6").Value = "第11章的阅读"
$ws.Range("H6").Value = "分章节、分层次、有自己理解"
$ws.Range("H6").WrapText = $true
$ws.Range("I6").Value = 4

# --- Row 7 (袁佳怡) -------------------------------------------------------
$ws.Range("B7").Value = 1
$ws.Range("D7").Value = 3
$ws.Range("F7").Value = "第11章的阅读"
$ws.Range("H7").Value = "有自己理解"
$ws.Range("I7").Value = 2

# --- Total-score formulas for rows 5-9 (shared group) --------------------
$ws.Range("L5").Formula = "=SUM(C5,E5,G5,I5,K5)"
$ws.Range("L6").Formula = "=SUM(C6,E6,G6,I6,K6)"
$ws.Range("L7").Formula = "=SUM(C7,E7,G7,I7,K7)"
$ws.Range("L8").Formula = "=SUM(C8,E8,G8,I8,K8)"
$ws.Range("L9").Formula = "=SUM(C9,E9,G9,I9,K9)"

# --- Move the active selection to J7, like the saved workbook ------------
[void]$ws.Range("J7").Select()
